# Insert a new weekly price record as row 439 (pushing all subsequent
# rows down by one) in the "Poroto verde" price sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 439; everything below shifts down.
$ws.Rows.Item(439).Insert()

# Populate the new row with the new record's data.
$ws.Cells.Item(439, 1).Value  = 3
$ws.Cells.Item(439, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(439, 3).Value  = "Coquimbo"
$ws.Cells.Item(439, 4).Value  = 44951
$ws.Cells.Item(439, 5).Value  = 5
$ws.Cells.Item(439, 6).Value  = 100112031
$ws.Cells.Item(439, 7).Value  = "Poroto verde"
$ws.Cells.Item(439, 8).Value  = "Magnum"
$ws.Cells.Item(439, 9).Value  = "Primera"
$ws.Cells.Item(439, 10).Value = 68
$ws.Cells.Item(439, 11).Value = 27000
$ws.Cells.Item(439, 12).Value = 28000
$ws.Cells.Item(439, 13).Value = 27559
$ws.Cells.Item(439, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(439, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(439, 16).Value = 1102
$ws.Cells.Item(439, 17).Value = 25
$ws.Cells.Item(439, 18).Value = "Hortaliza"
